# Fix AimType value names in EnemyAtkTable:
#   Angle -> World
#   PC    -> Target
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EnemyAtkGameData")

# Replace the mislabeled AimType values throughout the used range.
# xlWhole = 1, xlByRows = 1, xlNext = 1 (use literal constants since
# the Excel enum constants may not be available in this host).
$lookAt = 1      # xlWhole
$searchOrder = 1 # xlByRows
$searchDirection = 1 # xlNext

$ws.Cells.Replace("Angle", "World", $lookAt, $searchOrder, $false, $false, $false, $false)
$ws.Cells.Replace("PC", "Target", $lookAt, $searchOrder, $false, $false, $false, $false)

# Update the active selection to match the saved view state (G12).
$ws.Range("G12").Select()
